$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 4).Value2 = $ws.Cells.Item($row, 3).Value2
}
